$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 68

# Write the new run-log entry values first (on the un-styled row).
$ws.Cells.Item($row, 1).Value = "2025-08-28 09:38:37 UTC"
$ws.Cells.Item($row, 2).Value = "2025-08-28 15:08:37 IST"
$ws.Cells.Item($row, 3).Value = "SKIPPED"
$ws.Cells.Item($row, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($row, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-28-08-2025.pdf"
# Leading apostrophe forces a text cell whose displayed content is empty,
# matching the self-closing inlineStr cells (F/H) seen on the other
# SKIPPED rows.
$ws.Cells.Item($row, 6).Value = "'"
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = "'"

# Copy the formatting (style) of the previous data row onto the new row,
# same look as every other row in the log (centered alignment style).
$ws.Range("A67:H67").Copy()
$ws.Range("A68:H68").PasteSpecial(-4122)
